$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows right after the header (rows 2-4), pushing the existing
# 8 data rows down to rows 5-12. The sheets existing hyperlink refs (F2..F9,
# pointing at rId1..rId8) are left untouched by a native row insert, matching
# the target diff where those entries are unchanged.
$ws.Range("A2:A4").EntireRow.Insert()

# Row 2
$ws.Cells.Item(2,1).Value2 = '2026-01-30 18:38:56'
$ws.Cells.Item(2,2).Value2 = '製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)'
$ws.Cells.Item(2,3).Value2 = 'システム開発'
$ws.Cells.Item(2,4).Value2 = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(2,5).Value2 = '期限情報なし'
$ws.Cells.Item(2,6).Value2 = 'https://www.lancers.jp/work/detail/5460562'
$ws.Cells.Item(2,7).Value2 = 435
$ws.Cells.Item(2,8).Value2 = '🔥AI,Ai ◆ツール,開発'
$ws.Cells.Item(2,6).Style = "Hyperlink"

# Row 3
$ws.Cells.Item(3,1).Value2 = '2026-01-30 18:38:56'
$ws.Cells.Item(3,2).Value2 = '画像解析AI】釣具(ルアー)のオリジナルカラー判定システム開発(スマホ対応Web)'
$ws.Cells.Item(3,3).Value2 = 'システム開発'
$ws.Cells.Item(3,4).Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(3,5).Value2 = '期限情報なし'
$ws.Cells.Item(3,6).Value2 = 'https://www.lancers.jp/work/detail/5482607'
$ws.Cells.Item(3,7).Value2 = 410
$ws.Cells.Item(3,8).Value2 = '🔥AI,Ai ◆開発,システム開発'
$ws.Cells.Item(3,6).Style = "Hyperlink"

# Row 4
$ws.Cells.Item(4,1).Value2 = '2026-01-30 18:38:56'
$ws.Cells.Item(4,2).Value2 = '施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集'
$ws.Cells.Item(4,3).Value2 = 'システム開発'
$ws.Cells.Item(4,4).Value2 = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(4,5).Value2 = '期限情報なし'
$ws.Cells.Item(4,6).Value2 = 'https://www.lancers.jp/work/detail/5460563'
$ws.Cells.Item(4,7).Value2 = 220
$ws.Cells.Item(4,8).Value2 = '◆開発,システム開発 ◇管理'
$ws.Cells.Item(4,6).Style = "Hyperlink"

# Row 5
$ws.Cells.Item(5,1).Value2 = '2026-01-30 18:38:56'
$ws.Cells.Item(5,2).Value2 = '【Zapier保守・運用サポート】既存フローの管理・調整をお任せできる方募集(時給1,150円程度)'
$ws.Cells.Item(5,3).Value2 = 'システム開発'
$ws.Cells.Item(5,4).Value2 = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(5,5).Value2 = '期限情報なし'
$ws.Cells.Item(5,6).Value2 = 'https://www.lancers.jp/work/detail/5475245'
$ws.Cells.Item(5,7).Value2 = 213
$ws.Cells.Item(5,8).Value2 = '🔥API ◇管理'
$ws.Cells.Item(5,6).Style = "Hyperlink"

# Row 6
$ws.Cells.Item(6,1).Value2 = '2026-01-30 18:38:56'
$ws.Cells.Item(6,2).Value2 = 'シミュレーションスタジオの入退館システム開発(ロック選定含む/多店舗・複数打席対応)'
$ws.Cells.Item(6,3).Value2 = 'システム開発'
$ws.Cells.Item(6,4).Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(6,5).Value2 = '期限情報なし'
$ws.Cells.Item(6,6).Value2 = 'https://www.lancers.jp/work/detail/5482462'
$ws.Cells.Item(6,7).Value2 = 125
$ws.Cells.Item(6,8).Value2 = '◆開発,システム開発'
$ws.Cells.Item(6,6).Style = "Hyperlink"

# Row 7
$ws.Cells.Item(7,1).Value2 = '2026-01-30 18:38:56'
$ws.Cells.Item(7,2).Value2 = '【Java経験者】4月開始/与野 官公庁向けマイグレーション案件'
$ws.Cells.Item(7,3).Value2 = 'システム開発'
$ws.Cells.Item(7,4).Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(7,5).Value2 = '期限情報なし'
$ws.Cells.Item(7,6).Value2 = 'https://www.lancers.jp/work/detail/5482097'
$ws.Cells.Item(7,7).Value2 = 85
$ws.Cells.Item(7,8).Value2 = '★Java'
$ws.Cells.Item(7,6).Style = "Hyperlink"

# Row 8
$ws.Cells.Item(8,1).Value2 = '2026-01-30 18:38:56'
$ws.Cells.Item(8,2).Value2 = '店舗タブレット用Webサイトと管理画面の制作依頼'
$ws.Cells.Item(8,3).Value2 = 'システム開発'
$ws.Cells.Item(8,4).Value2 = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(8,5).Value2 = '期限情報なし'
$ws.Cells.Item(8,6).Value2 = 'https://www.lancers.jp/work/detail/5482389'
$ws.Cells.Item(8,7).Value2 = 63
$ws.Cells.Item(8,8).Value2 = '◇サイト'
$ws.Cells.Item(8,6).Style = "Hyperlink"

# Row 9
$ws.Cells.Item(9,1).Value2 = '2026-01-30 18:38:56'
$ws.Cells.Item(9,2).Value2 = 'WEBサーバーの管理、トラブル解決対応できる方を募集します!'
$ws.Cells.Item(9,3).Value2 = 'システム開発'
$ws.Cells.Item(9,4).Value2 = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(9,5).Value2 = '期限情報なし'
$ws.Cells.Item(9,6).Value2 = 'https://www.lancers.jp/work/detail/5481859'
$ws.Cells.Item(9,7).Value2 = 33
$ws.Cells.Item(9,8).Value2 = '◇管理'
$ws.Cells.Item(9,6).Style = "Hyperlink"

# Row 10
$ws.Cells.Item(10,1).Value2 = '2026-01-30 18:38:56'
$ws.Cells.Item(10,2).Value2 = '進行管理およびチームディレクションを担当'
$ws.Cells.Item(10,3).Value2 = 'システム開発'
$ws.Cells.Item(10,4).Value2 = '~ 5,000 円 / 固定'
$ws.Cells.Item(10,5).Value2 = '期限情報なし'
$ws.Cells.Item(10,6).Value2 = 'https://www.lancers.jp/work/detail/5418064'
$ws.Cells.Item(10,7).Value2 = 30
$ws.Cells.Item(10,8).Value2 = '◇管理'
$ws.Hyperlinks.Add($ws.Cells.Item(10,6), 'https://www.lancers.jp/work/detail/5418064')
$ws.Cells.Item(10,6).Style = "Hyperlink"

# Row 11
$ws.Cells.Item(11,1).Value2 = '2026-01-30 18:38:56'
$ws.Cells.Item(11,2).Value2 = '【高スキル】Web3.0系プロダクトの上級エンジニア募集'
$ws.Cells.Item(11,3).Value2 = 'システム開発'
$ws.Cells.Item(11,4).Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(11,5).Value2 = '期限情報なし'
$ws.Cells.Item(11,6).Value2 = 'https://www.lancers.jp/work/detail/5481715'
$ws.Cells.Item(11,7).Value2 = 25
$ws.Hyperlinks.Add($ws.Cells.Item(11,6), 'https://www.lancers.jp/work/detail/5481715')
$ws.Cells.Item(11,6).Style = "Hyperlink"

# Row 12
$ws.Cells.Item(12,1).Value2 = '2026-01-30 18:38:56'
$ws.Cells.Item(12,2).Value2 = '無人美容什器[ 決済 × IoT 連携の設計サポート ](※実装なし/スポット)'
$ws.Cells.Item(12,3).Value2 = 'システム開発'
$ws.Cells.Item(12,4).Value2 = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(12,5).Value2 = '期限情報なし'
$ws.Cells.Item(12,6).Value2 = 'https://www.lancers.jp/work/detail/5481888'
$ws.Cells.Item(12,7).Value2 = 18
$ws.Hyperlinks.Add($ws.Cells.Item(12,6), 'https://www.lancers.jp/work/detail/5481888')
$ws.Cells.Item(12,6).Style = "Hyperlink"

# Widen column H (skill summary) to fit the longer new values (XML width="19")
$ws.Columns.Item(8).ColumnWidth = 18.1666666667

